$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 18.79702060518519
$ws.Cells.Item(2, 3).Value = 9.357277552987691
$ws.Cells.Item(2, 4).Value = 6.004291376523517
$ws.Cells.Item(2, 5).Value = 10.04258830318587
$ws.Cells.Item(2, 7).Value = 44.88299066337746
$ws.Cells.Item(2, 8).Value = 17.80506267166255
$ws.Cells.Item(2, 9).Value = 26.92753150184702
$ws.Cells.Item(2, 12).Value = 10.26233251727938
$ws.Cells.Item(3, 2).Value = 18.27042250293245
$ws.Cells.Item(3, 3).Value = 8.798642518931992
$ws.Cells.Item(3, 4).Value = 5.890071163456551
$ws.Cells.Item(3, 5).Value = 10.06847501070571
$ws.Cells.Item(3, 7).Value = 44.50800301915827
$ws.Cells.Item(3, 8).Value = 17.8151023923383
$ws.Cells.Item(3, 9).Value = 26.97067776190817
$ws.Cells.Item(3, 12).Value = 10.23979208233152
$ws.Cells.Item(4, 2).Value = 17.94379365775438
$ws.Cells.Item(4, 3).Value = 8.435530188386023
$ws.Cells.Item(4, 4).Value = 5.820771634417095
$ws.Cells.Item(4, 5).Value = 10.08540831645012
$ws.Cells.Item(4, 7).Value = 44.29433128452798
$ws.Cells.Item(4, 8).Value = 17.8259964289664
$ws.Cells.Item(4, 9).Value = 27.00486733542242
$ws.Cells.Item(4, 12).Value = 10.2280645106213
$ws.Cells.Item(5, 2).Value = 17.81007815438944
$ws.Cells.Item(5, 3).Value = 8.282527720231753
$ws.Cells.Item(5, 4).Value = 5.792782341263018
$ws.Cells.Item(5, 5).Value = 10.09257036611906
$ws.Cells.Item(5, 7).Value = 44.21150149681134
$ws.Cells.Item(5, 8).Value = 17.8316205163199
$ws.Cells.Item(5, 9).Value = 27.02072409928497
$ws.Cells.Item(5, 12).Value = 10.22381939001256
$ws.Cells.Item(6, 2).Value = 17.78784442917313
$ws.Cells.Item(6, 3).Value = 8.256818809423564
$ws.Cells.Item(6, 4).Value = 5.788151104169619
$ws.Cells.Item(6, 5).Value = 10.09377542979803
$ws.Cells.Item(6, 7).Value = 44.19800596628582
$ws.Cells.Item(6, 8).Value = 17.83262579577386
$ws.Cells.Item(6, 9).Value = 27.0234729731628
$ws.Cells.Item(6, 12).Value = 10.2231468133556
$ws.Cells.Item(7, 2).Value = 17.94199250468943
$ws.Cells.Item(7, 3).Value = 8.433487075327289
$ws.Cells.Item(7, 4).Value = 5.820393090479529
$ws.Cells.Item(7, 5).Value = 10.08550384666888
$ws.Cells.Item(7, 7).Value = 44.29319694015868
$ws.Cells.Item(7, 8).Value = 17.82606748735831
$ws.Cells.Item(7, 9).Value = 27.00507340994318
$ws.Cells.Item(7, 12).Value = 10.2280050941168
$ws.Cells.Item(8, 2).Value = 18.61626008150546
$ws.Cells.Item(8, 3).Value = 9.16884611471629
$ws.Cells.Item(8, 4).Value = 5.964760456301621
$ws.Cells.Item(8, 5).Value = 10.05129867618786
$ws.Cells.Item(8, 7).Value = 44.75030409036246
$ws.Cells.Item(8, 8).Value = 17.80754014815012
$ws.Cells.Item(8, 9).Value = 26.94080408797524
$ws.Cells.Item(8, 12).Value = 10.2541239691651
$ws.Cells.Item(9, 2).Value = 19.90344368590331
$ws.Cells.Item(9, 3).Value = 10.45063307626951
$ws.Cells.Item(9, 4).Value = 6.25265216081317
$ws.Cells.Item(9, 5).Value = 9.992447160648329
$ws.Cells.Item(9, 7).Value = 45.77435256021106
$ws.Cells.Item(9, 8).Value = 17.80891445770826
$ws.Cells.Item(9, 9).Value = 26.87632236317329
$ws.Cells.Item(9, 12).Value = 10.32196151463589
$ws.Cells.Item(10, 2).Value = 20.81680350981256
$ws.Cells.Item(10, 3).Value = 11.29418003057665
$ws.Cells.Item(10, 4).Value = 6.464702499936552
$ws.Cells.Item(10, 5).Value = 9.954198859292001
$ws.Cells.Item(10, 7).Value = 46.59878187693874
$ws.Cells.Item(10, 8).Value = 17.83311502538452
$ws.Cells.Item(10, 9).Value = 26.86706312488953
$ws.Cells.Item(10, 12).Value = 10.3817150015959
$ws.Cells.Item(11, 2).Value = 21.22326111845021
$ws.Cells.Item(11, 3).Value = 11.65660391350637
$ws.Cells.Item(11, 4).Value = 6.560804253061339
$ws.Cells.Item(11, 5).Value = 9.937877201427371
$ws.Cells.Item(11, 7).Value = 46.98809545913478
$ws.Cells.Item(11, 8).Value = 17.84918907721053
$ws.Cells.Item(11, 9).Value = 26.87123460581167
$ws.Cells.Item(11, 12).Value = 10.41099529284691
$ws.Cells.Item(12, 2).Value = 21.37572717692166
$ws.Cells.Item(12, 3).Value = 11.79078286462864
$ws.Cells.Item(12, 4).Value = 6.597105547934759
$ws.Cells.Item(12, 5).Value = 9.931851220507282
$ws.Cells.Item(12, 7).Value = 47.13744605004199
$ws.Cells.Item(12, 8).Value = 17.85600531918581
$ws.Cells.Item(12, 9).Value = 26.87402687707463
$ws.Cells.Item(12, 12).Value = 10.42237925375544
$ws.Cells.Item(13, 2).Value = 21.34295766080341
$ws.Cells.Item(13, 3).Value = 11.76202121221009
$ws.Cells.Item(13, 4).Value = 6.589292024950765
$ws.Cells.Item(13, 5).Value = 9.933142148629139
$ws.Cells.Item(13, 7).Value = 47.10519713813672
$ws.Cells.Item(13, 8).Value = 17.85450487396097
$ws.Cells.Item(13, 9).Value = 26.87337147852489
$ws.Cells.Item(13, 12).Value = 10.41991443562365
$ws.Cells.Item(14, 2).Value = 21.23583448250066
$ws.Cells.Item(14, 3).Value = 11.66770426952293
$ws.Cells.Item(14, 4).Value = 6.563792806705206
$ws.Cells.Item(14, 5).Value = 9.93737834205451
$ws.Cells.Item(14, 7).Value = 47.00034476568958
$ws.Cells.Item(14, 8).Value = 17.84973523339185
$ws.Cells.Item(14, 9).Value = 26.87143998818945
$ws.Cells.Item(14, 12).Value = 10.41192595617859
$ws.Cells.Item(15, 2).Value = 21.17002527672738
$ws.Cells.Item(15, 3).Value = 11.60953359295934
$ws.Cells.Item(15, 4).Value = 6.548160923915132
$ws.Cells.Item(15, 5).Value = 9.939993267729067
$ws.Cells.Item(15, 7).Value = 46.93636663469187
$ws.Cells.Item(15, 8).Value = 17.84690868968113
$ws.Cells.Item(15, 9).Value = 26.87041500546374
$ws.Cells.Item(15, 12).Value = 10.40707117731037
$ws.Cells.Item(16, 2).Value = 20.79004547892821
$ws.Cells.Item(16, 3).Value = 11.27006601169984
$ws.Cells.Item(16, 4).Value = 6.458411431404411
$ws.Cells.Item(16, 5).Value = 9.955287178594457
$ws.Cells.Item(16, 7).Value = 46.57361585367653
$ws.Cells.Item(16, 8).Value = 17.83216657339867
$ws.Cells.Item(16, 9).Value = 26.86695983671745
$ws.Cells.Item(16, 12).Value = 10.37984321980903
$ws.Cells.Item(17, 2).Value = 20.55451393074309
$ws.Cells.Item(17, 3).Value = 11.05635567635702
$ws.Cells.Item(17, 4).Value = 6.403232626592005
$ws.Cells.Item(17, 5).Value = 9.964945307987067
$ws.Cells.Item(17, 7).Value = 46.3546477690637
$ws.Cells.Item(17, 8).Value = 17.82442100517164
$ws.Cells.Item(17, 9).Value = 26.8669928054689
$ws.Cells.Item(17, 12).Value = 10.36367330951935
$ws.Cells.Item(18, 2).Value = 20.41820019289264
$ws.Cells.Item(18, 3).Value = 10.93142909289487
$ws.Cells.Item(18, 4).Value = 6.371463445320893
$ws.Cells.Item(18, 5).Value = 9.970601862630305
$ws.Cells.Item(18, 7).Value = 46.23005483521194
$ws.Cells.Item(18, 8).Value = 17.82044269914052
$ws.Cells.Item(18, 9).Value = 26.86780051193898
$ws.Cells.Item(18, 12).Value = 10.35457065058314
$ws.Cells.Item(19, 2).Value = 20.37190690157765
$ws.Cells.Item(19, 3).Value = 10.88878656457537
$ws.Cells.Item(19, 4).Value = 6.360702687702775
$ws.Cells.Item(19, 5).Value = 9.972534508713494
$ws.Cells.Item(19, 7).Value = 46.18810575735133
$ws.Cells.Item(19, 8).Value = 17.81917755065855
$ws.Cells.Item(19, 9).Value = 26.86820922173289
$ws.Cells.Item(19, 12).Value = 10.35152278937593
$ws.Cells.Item(20, 2).Value = 20.5796749754851
$ws.Cells.Item(20, 3).Value = 11.07931312458795
$ws.Cells.Item(20, 4).Value = 6.409110079330377
$ws.Cells.Item(20, 5).Value = 9.96390668566416
$ws.Cells.Item(20, 7).Value = 46.37781830286102
$ws.Cells.Item(20, 8).Value = 17.8251961837295
$ws.Cells.Item(20, 9).Value = 26.86690760927181
$ws.Cells.Item(20, 12).Value = 10.36537418673414
$ws.Cells.Item(21, 2).Value = 21.26733966274397
$ws.Cells.Item(21, 3).Value = 11.69549055103409
$ws.Cells.Item(21, 4).Value = 6.571285294347226
$ws.Cells.Item(21, 5).Value = 9.936129873781606
$ws.Cells.Item(21, 7).Value = 47.03109121249808
$ws.Cells.Item(21, 8).Value = 17.8511163940686
$ws.Cells.Item(21, 9).Value = 26.87197435242658
$ws.Cells.Item(21, 12).Value = 10.41426437515157
$ws.Cells.Item(22, 2).Value = 21.70825264988551
$ws.Cells.Item(22, 3).Value = 12.08034629106341
$ws.Cells.Item(22, 4).Value = 6.676733061592759
$ws.Cells.Item(22, 5).Value = 9.918877564378137
$ws.Cells.Item(22, 7).Value = 47.46920713869386
$ws.Cells.Item(22, 8).Value = 17.8723078597567
$ws.Cells.Item(22, 9).Value = 26.88235636965407
$ws.Cells.Item(22, 12).Value = 10.44794026379009
$ws.Cells.Item(23, 2).Value = 21.4737546661529
$ws.Cells.Item(23, 3).Value = 11.87657392640429
$ws.Cells.Item(23, 4).Value = 6.620515501637493
$ws.Cells.Item(23, 5).Value = 9.928003059617046
$ws.Cells.Item(23, 7).Value = 47.23439793585724
$ws.Cells.Item(23, 8).Value = 17.86060847017551
$ws.Cells.Item(23, 9).Value = 26.87616627122978
$ws.Cells.Item(23, 12).Value = 10.42981106584231
$ws.Cells.Item(24, 2).Value = 20.56830246975201
$ws.Cells.Item(24, 3).Value = 11.06894048393372
$ws.Cells.Item(24, 4).Value = 6.406453023761191
$ws.Cells.Item(24, 5).Value = 9.96437592295862
$ws.Cells.Item(24, 7).Value = 46.36733886214562
$ws.Cells.Item(24, 8).Value = 17.82484424692748
$ws.Cells.Item(24, 9).Value = 26.866943670234
$ws.Cells.Item(24, 12).Value = 10.3646046165901
$ws.Cells.Item(25, 2).Value = 19.56016632932303
$ws.Cells.Item(25, 3).Value = 10.12110169025678
$ws.Cells.Item(25, 4).Value = 6.174507204825712
$ws.Cells.Item(25, 5).Value = 10.00749005820099
$ws.Cells.Item(25, 7).Value = 45.48421609408025
$ws.Cells.Item(25, 8).Value = 17.80448205327568
$ws.Cells.Item(25, 9).Value = 26.88711210019193
$ws.Cells.Item(25, 12).Value = 10.30185120654788
